$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.235148
$ws.Range("H2").Value = 21.705444
$ws.Range("I2").Value = 0.9254344869740032
$ws.Range("J2").Value = 0.9254344869740032
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.380049
$ws.Range("N2").Value = 100.140147
$ws.Range("O2").Value = 0.3891462059670435
$ws.Range("P2").Value = 0.3891462059670435
$ws.Range("Q2").Value = 241.509594762252
$ws.Range("R2").Value = 2173.586352860268
$ws.Range("S2").Value = 0.3601293194769907
$ws.Range("T2").Value = 0.3601293194769907

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.235148
$ws.Range("H3").Value = 21.705444
$ws.Range("I3").Value = 0.9254344869740032
$ws.Range("J3").Value = 0.9254344869740032
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.85786133333334
$ws.Range("N3").Value = 101.573584
$ws.Range("O3").Value = 0.3947165649764305
$ws.Range("P3").Value = 0.3947165649764305
$ws.Range("Q3").Value = 244.966637710144
$ws.Range("R3").Value = 2204.699739391296
$ws.Range("S3").Value = 0.3652843218091038
$ws.Range("T3").Value = 0.3652843218091038

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.235148
$ws.Range("H4").Value = 21.705444
$ws.Range("I4").Value = 0.9254344869740032
$ws.Range("J4").Value = 0.9254344869740032
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.53974466666667
$ws.Range("N4").Value = 55.61923400000001
$ws.Range("O4").Value = 0.2161372290565261
$ws.Range("P4").Value = 0.2161372290565261
$ws.Range("Q4").Value = 134.137796545544
$ws.Range("R4").Value = 1207.240168909896
$ws.Range("S4").Value = 0.2000208456879088
$ws.Range("T4").Value = 0.2000208456879088

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3016356666666667
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.03858166390441884
$ws.Range("J5").Value = 0.03858166390441884
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.380049
$ws.Range("N5").Value = 100.140147
$ws.Range("O5").Value = 0.3891462059670435
$ws.Range("P5").Value = 0.3891462059670435
$ws.Range("Q5").Value = 10.068613333481
$ws.Range("R5").Value = 90.61752000132901
$ws.Range("S5").Value = 0.01501390812830022
$ws.Range("T5").Value = 0.01501390812830022

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3016356666666667
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.03858166390441884
$ws.Range("J6").Value = 0.03858166390441884
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 33.85786133333334
$ws.Range("N6").Value = 101.573584
$ws.Range("O6").Value = 0.3947165649764305
$ws.Range("P6").Value = 0.3947165649764305
$ws.Range("Q6").Value = 10.21273857518756
$ws.Range("R6").Value = 91.91464717668802
$ws.Range("S6").Value = 0.01522882184742734
$ws.Range("T6").Value = 0.01522882184742734

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3016356666666667
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.03858166390441884
$ws.Range("J7").Value = 0.03858166390441884
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.53974466666667
$ws.Range("N7").Value = 55.61923400000001
$ws.Range("O7").Value = 0.2161372290565261
$ws.Range("P7").Value = 0.2161372290565261
$ws.Range("Q7").Value = 5.592248242359778
$ws.Range("R7").Value = 50.33023418123801
$ws.Range("S7").Value = 0.00833893392869128
$ws.Range("T7").Value = 0.00833893392869128

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.0359838491215779
$ws.Range("J8").Value = 0.0359838491215779
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.380049
$ws.Range("N8").Value = 100.140147
$ws.Range("O8").Value = 0.3891462059670435
$ws.Range("P8").Value = 0.3891462059670435
$ws.Range("Q8").Value = 9.390664538290999
$ws.Range("R8").Value = 84.515980844619
$ws.Range("S8").Value = 0.01400297836175257
$ws.Range("T8").Value = 0.01400297836175257

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.0359838491215779
$ws.Range("J9").Value = 0.0359838491215779
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.85786133333334
$ws.Range("N9").Value = 101.573584
$ws.Range("O9").Value = 0.3947165649764305
$ws.Range("P9").Value = 0.3947165649764305
$ws.Range("Q9").Value = 9.525085411507556
$ws.Range("R9").Value = 85.72576870356801
$ws.Range("S9").Value = 0.01420342131989938
$ws.Range("T9").Value = 0.01420342131989937

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.0359838491215779
$ws.Range("J10").Value = 0.0359838491215779
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.53974466666667
$ws.Range("N10").Value = 55.61923400000001
$ws.Range("O10").Value = 0.2161372290565261
$ws.Range("P10").Value = 0.2161372290565261
$ws.Range("Q10").Value = 5.215706028179778
$ws.Range("R10").Value = 46.941354253618
$ws.Range("S10").Value = 0.007777449439925957
$ws.Range("T10").Value = 0.007777449439925957
